# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
# Mirrors the commit "Created functions to get season record": the sheet
# previously only had team/player stats through column AC; this appends
# three new columns (AD:AF) with the team's season W-L-T record repeated
# on every player row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) -------------------------------------------------
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the bold/centered/bordered header style already used by the rest
# of row 1 (e.g. AC1) by copying its format onto the new header cells.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# --- Data rows (rows 2-42): team's 2006 season record = 89-73-0 --------
$lastRow = 42
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 89   # AD: Wins
    $ws.Cells.Item($r, 31).Value = 73   # AE: Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF: Ties
}
